$d = $word.ActiveDocument

$replacements = @(
    @('77×76=', '51×48='),
    @('24×17=', '49×12='),
    @('49×87=', '31×44='),
    @('37×27=', '55×46='),
    @('71×55=', '78×43='),
    @('35×85=', '22×16='),
    @('44×99=', '89×32='),
    @('28×85=', '55×16='),
    @('98×34=', '19×61='),
    @('54×84=', '96×13='),
    @('32×81=', '61×11='),
    @('42×72=', '93×40='),
    @('72×40=', '97×44='),
    @('54×69=', '56×17='),
    @('66×25=', '34×29='),
    @('63×37=', '83×80='),
    @('42×76=', '84×54='),
    @('82×64=', '93×63='),
    @('73×32=', '45×96='),
    @('79×42=', '45×67='),
    @('52×68=', '26×68='),
    @('14×67=', '85×44='),
    @('40×81=', '37×63='),
    @('55×68=', '72×46='),
    @('82×60=', '88×86=')
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

$d.Save()
